$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the text of the existing entry describing 12.11. (row 15 / shared string 16) ---
$ws.Cells.Item(15, 3).Value = "Purkualgoritmin koodausta, debuggausta. Purkualgoritmi mahdollisesti toimiikin jo oikein mutta Huffman-koodin pakkaus tavujen biteiksi lienee buginen. Nyt sekä pakkaus että purku toimii melkein: purussa vielä bugi joka tuottaa epämääräisiä virheitä striimin sekaan. Ongelma liittyy Huffman-koodiin joka luetaan kahdessa lohkossa levyltä (todiste: kun lohkokoko >= pakatun tiedoston koko, ongelma häviää)."

# Row 15 (12.11.2012): hours 5 -> 11, taller row to fit the longer text
$ws.Cells.Item(15, 2).Value = 11
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 (13.11.2012): new entry ---
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 41226
$ws.Cells.Item(16, 2).Value = 8
$ws.Cells.Item(16, 3).Value = "Uudelleenkirjoitettu osia pakkausalgoritmista joka korjaa eilisen bugin. Nyt sekä pakkaus että purku toimii 100%. Hajautustaulukko kirjoitettu (iteraattorit avaimille ja arvoille puuttuu)."
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17 (14.11.2012): new entry ---
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).Value = 41227
$ws.Cells.Item(17, 2).Value = 3
$ws.Cells.Item(17, 3).Value = "Hajautustaulukon iteraattorit toteutettu ja testattu.`nKirjoitettu  yksikkötesti pääluokalle. Testattu erikokoisilla syötteillä,  käytetty testidata lisätty myös GitHubiin. JavaDoceja viimeistelty. Jäljellä työssä  prioriteettijonon toteutus."
$ws.Rows.Item(17).RowHeight = 75

$excel.CutCopyMode = $false

# B24 (Yhteensä) recalculates automatically via its existing SUM(B6:B23) formula.

# --- Update the view state: scrolled down, C18 selected ---
$ws.Activate()
$excel.Goto($ws.Range("A14"), $true)
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18").Select()
